# Append the new allocation row produced by the 2025-09-08 run to the
# bottom of the table (the sheet currently ends at row 6, data goes to row 7).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A stores the date as a literal "MM/DD/YYYY" text label (same as
# every row above it), not a real Excel date value/serial. A leading
# apostrophe forces Excel to keep the value as literal text instead of
# auto-converting the date-like string into a date serial number.
# ClearFormats() then drops the "quote prefix" cell formatting that the
# apostrophe trick applies, so the cell ends up unstyled - just like the
# existing data rows.
$ws.Range("A7").Value = "'09/08/2025"
$ws.Range("A7").ClearFormats()

$ws.Range("B7").Value = 0.1229112743299529
$ws.Range("C7").Value = 0.8770887256700471
